$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for column C (y_1_forecast) and column E (y_1_forecast second series)
# keyed by row number, rows 4..19
$values = @{
    4  = @{ C = 3.184002331674129;   E = 6.704254199558113 }
    5  = @{ C = 8.626810748872327;   E = 4.739201070534826 }
    6  = @{ C = 3.449685446853534;   E = 3.26507595662513 }
    7  = @{ C = -1.480934717826909;  E = 0.7772706050320544 }
    8  = @{ C = 5.427992542801308;   E = 4.945882057432871 }
    9  = @{ C = 4.970284184513551;   E = 4.488174889976171 }
    10 = @{ C = 4.039484738713828;   E = 4.214976960249173 }
    11 = @{ C = 4.589070866863865;   E = 3.829046580278361 }
    12 = @{ C = 3.625873842174787;   E = 2.330842103296149 }
    13 = @{ C = 3.500574054404404;   E = 3.21661481720994 }
    14 = @{ C = -0.9913189363815245; E = 1.183532150252908 }
    15 = @{ C = 4.507091823899212;   E = 5.429743376942153 }
    16 = @{ C = 8.053468068361846;   E = 3.974997080343634 }
    17 = @{ C = 0.2714278794373248;  E = 2.460471645027118 }
    18 = @{ C = -2.107534670984712;  E = 2.747596279389564 }
    19 = @{ C = 1.552685227480533;   E = 2.496145622272206 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
